$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD, AE, AF) after the existing last column (AC),
# copying the header style (s="1") from the existing header cell AC1 so no new
# style entries are introduced.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row (2-49).
$ws.Range("AD2:AD49").Value = 77
$ws.Range("AE2:AE49").Value = 85
$ws.Range("AF2:AF49").Value = 0
